$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of the last existing data row (row 28) down onto the
# six new rows (29-34) so the new cells pick up the same cell styles
# (centered "Time in hrs" column, date-formatted "Due Date" column, etc.)
# that the rest of the table uses.
$ws.Range("B28:D28").Copy() | Out-Null
$ws.Range("B29:D34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Due dates were filled in first, straight down the column
$ws.Range("D29").Value = " 5/17/2020"
$ws.Range("D30").Value = " 5/17/2021"
$ws.Range("D31").Value = " 5/17/2022"
$ws.Range("D32").Value = " 5/17/2023"
$ws.Range("D33").Value = " 5/17/2024"
$ws.Range("D34").Value = " 5/17/2025"

# New assignment rows (name + hours) filled in afterwards
$ws.Range("B29").Value = "MVC Core SportsStore App, 2 - Chapter 10"
$ws.Range("C29").Value = 2

$ws.Range("B30").Value = "MVC Core SportsStore App, 2 - Chapter 11"
$ws.Range("C30").Value = 3

# Grow the table / autofilter to cover the newly added rows
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B2:F34"))

# Restore view state: scroll the window down a bit and leave the active
# selection on C31 (mirrors scrolling to row 13 and clicking C31).
$ws.Activate() | Out-Null
$ws.Range("C31").Select() | Out-Null
try { $excel.ActiveWindow.ScrollRow = 13 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
